$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the FlashScore odds data refresh (commit: "Atualizando o arquivo XLSX").
# Each assignment sets <col><row> to its new numeric odd/value as specified by the diff.

# Row 7
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 15

# Row 13
$ws.Range("AC13").Value = 21
$ws.Range("AD13").Value = 8.5
$ws.Range("AE13").Value = 11
$ws.Range("AI13").Value = 13
$ws.Range("AK13").Value = 19
$ws.Range("AN13").Value = 6
$ws.Range("AO13").Value = 17
$ws.Range("AP13").Value = 19
$ws.Range("AX13").Value = 9.5
$ws.Range("G13").Value = 3.4
$ws.Range("H13").Value = 4.1
$ws.Range("I13").Value = 1.91
$ws.Range("J13").Value = 3.6
$ws.Range("L13").Value = 2.4
$ws.Range("U13").Value = 1.4
$ws.Range("V13").Value = 2.75

# Row 15
$ws.Range("BA15").Value = 67
$ws.Range("I15").Value = 3.7
$ws.Range("L15").Value = 4
$ws.Range("Q15").Value = 1.53
$ws.Range("R15").Value = 2.4
$ws.Range("U15").Value = 1.5
$ws.Range("V15").Value = 2.5
$ws.Range("W15").Value = 11
$ws.Range("X15").Value = 11

# Row 16
$ws.Range("AG16").Value = 101
$ws.Range("AI16").Value = 17
$ws.Range("AK16").Value = 29
$ws.Range("AM16").Value = 23
$ws.Range("AN16").Value = 4.75
$ws.Range("AO16").Value = 12
$ws.Range("AR16").Value = 41
$ws.Range("AY16").Value = 19
$ws.Range("AZ16").Value = 41
$ws.Range("BB16").Value = 101
$ws.Range("BC16").Value = 301
$ws.Range("G16").Value = 2.2
$ws.Range("H16").Value = 3.75
$ws.Range("I16").Value = 2.88
$ws.Range("K16").Value = 2.4
$ws.Range("L16").Value = 3.25
$ws.Range("Q16").Value = 1.53
$ws.Range("R16").Value = 2.4
$ws.Range("U16").Value = 1.5
$ws.Range("V16").Value = 2.5
$ws.Range("W16").Value = 12
$ws.Range("X16").Value = 13
$ws.Range("Y16").Value = 9.5

# Row 17
$ws.Range("AC17").Value = 12
$ws.Range("AE17").Value = 17
$ws.Range("AG17").Value = 301
$ws.Range("AK17").Value = 51
$ws.Range("AO17").Value = 7.5
$ws.Range("AS17").Value = 126
$ws.Range("AT17").Value = 3.25
$ws.Range("BD17").Value = 151
$ws.Range("I17").Value = 5.5
$ws.Range("O17").Value = 1.22
$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 1.75
$ws.Range("R17").Value = 2.05
$ws.Range("S17").Value = 1.33
$ws.Range("T17").Value = 3.25
$ws.Range("U17").Value = 1.83
$ws.Range("V17").Value = 1.83
$ws.Range("W17").Value = 7
$ws.Range("X17").Value = 7.5

# Row 18
$ws.Range("AD18").Value = 8.5
$ws.Range("AH18").Value = 17
$ws.Range("AI18").Value = 21
$ws.Range("AJ18").Value = 12
$ws.Range("AK18").Value = 34
$ws.Range("AL18").Value = 21
$ws.Range("AP18").Value = 15
$ws.Range("AT18").Value = 4
$ws.Range("AV18").Value = 34
$ws.Range("G18").Value = 2.1
$ws.Range("H18").Value = 3.9
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 2.6
$ws.Range("L18").Value = 3.25
$ws.Range("S18").Value = 1.22
$ws.Range("T18").Value = 4
$ws.Range("Z18").Value = 21

# Row 19
$ws.Range("AA19").Value = 17
$ws.Range("AC19").Value = 15
$ws.Range("AH19").Value = 12
$ws.Range("AI19").Value = 17
$ws.Range("AJ19").Value = 11
$ws.Range("AN19").Value = 4.5
$ws.Range("AO19").Value = 12
$ws.Range("AP19").Value = 19
$ws.Range("AS19").Value = 101
$ws.Range("AT19").Value = 3.4
$ws.Range("AU19").Value = 7
$ws.Range("BD19").Value = 151
$ws.Range("G19").Value = 2.3
$ws.Range("I19").Value = 2.8
$ws.Range("J19").Value = 2.88
$ws.Range("M19").Value = 1.03
$ws.Range("N19").Value = 15
$ws.Range("O19").Value = 1.18
$ws.Range("P19").Value = 4.5
$ws.Range("Q19").Value = 1.65
$ws.Range("R19").Value = 2.2
$ws.Range("S19").Value = 1.3
$ws.Range("T19").Value = 3.4
$ws.Range("W19").Value = 11

# Row 20
$ws.Range("AG20").Value = 201
$ws.Range("AH20").Value = 15
$ws.Range("AK20").Value = 51
$ws.Range("AN20").Value = 3.75
$ws.Range("AO20").Value = 8.5
$ws.Range("AQ20").Value = 26
$ws.Range("AU20").Value = 8
$ws.Range("AY20").Value = 29
$ws.Range("AZ20").Value = 81
$ws.Range("G20").Value = 1.7
$ws.Range("I20").Value = 4.33
$ws.Range("L20").Value = 4.75
$ws.Range("U20").Value = 1.73
$ws.Range("V20").Value = 2
$ws.Range("W20").Value = 8
$ws.Range("Z20").Value = 13

# Row 21
$ws.Range("AB21").Value = 29
$ws.Range("AC21").Value = 21
$ws.Range("AD21").Value = 9.5
$ws.Range("AF21").Value = 34
$ws.Range("AG21").Value = 101
$ws.Range("AL21").Value = 11
$ws.Range("AR21").Value = 67
$ws.Range("AT21").Value = 4
$ws.Range("BA21").Value = 34
$ws.Range("H21").Value = 4.5
$ws.Range("K21").Value = 2.63
$ws.Range("L21").Value = 2.05
$ws.Range("N21").Value = 21
$ws.Range("O21").Value = 1.13
$ws.Range("P21").Value = 6
$ws.Range("Q21").Value = 1.44
$ws.Range("R21").Value = 2.7
$ws.Range("S21").Value = 1.22
$ws.Range("T21").Value = 4
$ws.Range("U21").Value = 1.5
$ws.Range("V21").Value = 2.5
$ws.Range("W21").Value = 21

# Row 22
$ws.Range("AC22").Value = 26
$ws.Range("AD22").Value = 13
$ws.Range("AE22").Value = 19
$ws.Range("AG22").Value = 151
$ws.Range("AJ22").Value = 23
$ws.Range("AL22").Value = 51
$ws.Range("AQ22").Value = 13
$ws.Range("AW22").Value = 9.5
$ws.Range("AZ22").Value = 126
$ws.Range("G22").Value = 1.29
$ws.Range("H22").Value = 6
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 1.67
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 7
$ws.Range("N22").Value = 26
$ws.Range("U22").Value = 1.62
$ws.Range("V22").Value = 2.2

# Row 25
$ws.Range("Q25").Value = 1.93
$ws.Range("R25").Value = 1.93

# Row 26
$ws.Range("AB26").Value = 26
$ws.Range("AM26").Value = 29
$ws.Range("AT26").Value = 3
$ws.Range("AV26").Value = 51
$ws.Range("BB26").Value = 151
$ws.Range("G26").Value = 2.38
$ws.Range("N26").Value = 12
$ws.Range("S26").Value = 1.36
$ws.Range("T26").Value = 3

# Row 27
$ws.Range("AA27").Value = 19
$ws.Range("AP27").Value = 26
$ws.Range("O27").Value = 1.44
$ws.Range("P27").Value = 2.63

# Row 30
$ws.Range("Q30").Value = 1.93
$ws.Range("R30").Value = 1.93

# Row 31
$ws.Range("BD31").Value = 176
$ws.Range("O31").Value = 1.14
$ws.Range("P31").Value = 5.5

# Row 33
$ws.Range("O33").Value = 1.2
$ws.Range("P33").Value = 4.33
$ws.Range("Q33").Value = 1.67
$ws.Range("R33").Value = 2.15
